# Update For PT Setting
# Adds two new worksheets (CAPTSetting, SMAPTSetting), appends a new product
# ("Play n Go") to PositiveExtra / ProductList, and populates the new sheets
# with per-product PT-setting data.

$wb = $excel.ActiveWorkbook
$wsPositive = $wb.Worksheets.Item(1)   # PositiveExtra
$wsProduct  = $wb.Worksheets.Item(2)   # ProductList

# ---------------------------------------------------------------------------
# 1. PositiveExtra: append row 16 (A16 = 15, B16 = 10.5)
# ---------------------------------------------------------------------------
$wsPositive.Range("A16").NumberFormat = $wsPositive.Range("A15").NumberFormat
$wsPositive.Range("A16").Value = 15
$wsPositive.Range("B16").Value = 10.5

# ---------------------------------------------------------------------------
# 2. Create the two new worksheets, right after ProductList, in order:
#    CAPTSetting, then SMAPTSetting
# ---------------------------------------------------------------------------
$wsCAPT = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsProduct)
$wsCAPT.Name = "CAPTSetting"

$wsSMAPT = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCAPT)
$wsSMAPT.Name = "SMAPTSetting"

# Header look (bold + yellow fill, same as existing headers) for both sheets.
foreach ($ws in @($wsCAPT, $wsSMAPT)) {
    $ws.Range("A1:E1").Font.Bold = $true
    $ws.Range("A1:E1").Interior.Color = $wsProduct.Range("A1").Interior.Color
}

# ---------------------------------------------------------------------------
# 3. Header text. New shared strings get interned in the exact order they
#    are first assigned, so we sequence these writes (interleaved with step
#    4) to land them at shared-string indices 17 (PPT1), 18 (Play n Go),
#    19 (NPT1), 20 (NPT2) -- matching the target file.
# ---------------------------------------------------------------------------
$wsCAPT.Range("A1").Value = "ProdID"   # reuses existing shared string 0
$wsCAPT.Range("B1").Value = "Name"     # reuses existing shared string 2
$wsCAPT.Range("C1").Value = "PPT1"     # -> new shared string 17

# ---------------------------------------------------------------------------
# 4. ProductList: append row 16 (A16 = 15, B16 = "Play n Go")
#    ("Play n Go" becomes shared-string index 18.)
# ---------------------------------------------------------------------------
$wsProduct.Range("A16").NumberFormat = $wsProduct.Range("A15").NumberFormat
$wsProduct.Range("A16").Value = 15
$wsProduct.Range("B16").Value = "Play n Go"

# ---------------------------------------------------------------------------
# 5. Finish the CAPTSetting / SMAPTSetting header rows.
# ---------------------------------------------------------------------------
$wsCAPT.Range("D1").Value = "NPT1"     # -> new shared string 19
$wsCAPT.Range("E1").Value = "NPT2"     # -> new shared string 20

$wsSMAPT.Range("A1").Value = "ProdID"
$wsSMAPT.Range("B1").Value = "Name"
$wsSMAPT.Range("C1").Value = "PPT1"
$wsSMAPT.Range("D1").Value = "NPT1"
$wsSMAPT.Range("E1").Value = "NPT2"

# ---------------------------------------------------------------------------
# 6. Populate CAPTSetting / SMAPTSetting data rows (2-16), column A/B mirror
#    ProductList, columns C/D/E hold the fixed PT-setting values per sheet.
# ---------------------------------------------------------------------------
$names = @()
for ($r = 2; $r -le 16; $r++) {
    $names += $wsProduct.Range("B$r").Value2
}

for ($i = 0; $i -lt 15; $i++) {
    $r = $i + 2
    $name = $names[$i]

    $wsCAPT.Range("A$r").NumberFormat = $wsProduct.Range("A2").NumberFormat
    $wsCAPT.Range("A$r").Value = $i + 1
    $wsCAPT.Range("B$r").Value = $name
    $wsCAPT.Range("C$r").Value = 50.5
    $wsCAPT.Range("D$r").Value = 49
    $wsCAPT.Range("E$r").Value = 100

    $wsSMAPT.Range("A$r").NumberFormat = $wsProduct.Range("A2").NumberFormat
    $wsSMAPT.Range("A$r").Value = $i + 1
    $wsSMAPT.Range("B$r").Value = $name
    $wsSMAPT.Range("C$r").Value = 20
    $wsSMAPT.Range("D$r").Value = 19
    $wsSMAPT.Range("E$r").Value = 80
}

# ---------------------------------------------------------------------------
# 7. Column widths on the new sheets' Name column (best-effort; matches the
#    custom widths used in the target worksheet).
# ---------------------------------------------------------------------------
$wsCAPT.Columns.Item(2).ColumnWidth = 17.11
$wsSMAPT.Columns.Item(2).ColumnWidth = 16.55

# ---------------------------------------------------------------------------
# 8. Selections / active sheet to match the target view state as closely as
#    the host allows.
# ---------------------------------------------------------------------------
$wsCAPT.Activate()
$wsCAPT.Range("E2").Select()

$wsSMAPT.Activate()
$wsSMAPT.Range("F13").Select()

$wsProduct.Activate()
$wsProduct.Range("A16:B16").Select()

$wsPositive.Activate()
$wsPositive.Range("B17").Select()
